$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeById($shape, $id) {
    if ($shape.Id -eq $id) {
        return $shape
    }
    if ($shape.Type -eq 6) {
        $items = $shape.GroupItems
        for ($i = 1; $i -le $items.Count; $i++) {
            $found = Find-ShapeById $items.Item($i) $id
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

function Get-ShapeById($id) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $found = Find-ShapeById $s.Shapes.Item($i) $id
        if ($found -ne $null) {
            return $found
        }
    }
    return $null
}

# The "TCPCEP_001" boxes (second repetition box, both in the "with Dynamic
# Connection" and "without Dynamic Connection" groups) incorrectly repeat the
# REP_000 / REP_001 labels used by the first "TCPCEP_000" box. Fix the
# figures so the second box shows REP_002 / REP_003.
$targets = @(
    @{ Id = 33; Old = "REP_000"; New = "REP_002" },
    @{ Id = 31; Old = "REP_001"; New = "REP_003" },
    @{ Id = 87; Old = "REP_000"; New = "REP_002" },
    @{ Id = 89; Old = "REP_001"; New = "REP_003" }
)

foreach ($t in $targets) {
    $shape = Get-ShapeById $t.Id
    if ($shape -eq $null) { continue }
    if (-not $shape.HasTextFrame) { continue }
    $tr = $shape.TextFrame.TextRange
    for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
        $para = $tr.Paragraphs($pi)
        if ($para.Text -eq $t.Old) {
            $para.Runs(1).Text = $t.New
        }
    }
}
